$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Overview" sheet: update the rollup status text for both locale columns
#    ("Ready for handoff" -> "Handed back: in sync with en-US").
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$newStatus = "Handed back: in sync with en-US"
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# Common URLs reused by both locale sheets' hyperlinks.
# ---------------------------------------------------------------------------
$urlMd = "https://github.com/OpenLocalizationTest/oltest/blob/ca4cd978eab06060218aac716cd58ab6e24090fd/e2e/28a0734b-b7aa-416f-b69c-8b1b21830f38.md"
$urlFfff = "https://github.com/OpenLocalizationTest/oltest/blob/ca4cd978eab06060218aac716cd58ab6e24090fd/e2e/ffff97bd2b86-d13f-492f-964b-72d34ec1936d.md"
$urlConfig = "https://github.com/OpenLocalizationTest/oltest/blob/ca4cd978eab06060218aac716cd58ab6e24090fd/.localization-config"
$urlZhXlf = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fa218e84d0a5e51ad3cdbe9d62d8481294d25ea8/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/28a0734b-b7aa-416f-b69c-8b1b21830f38.cb1e54aa0bb61f184d790a69febda753cb09c03f.zh-cn.xlf"
$urlDeXlf = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d3df459718bbfaeeeba3e9988074075382f6a55f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/28a0734b-b7aa-416f-b69c-8b1b21830f38.cb1e54aa0bb61f184d790a69febda753cb09c03f.de-de.xlf"

$mdName = "28a0734b-b7aa-416f-b69c-8b1b21830f38.md"
$ffffName = "ffff97bd2b86-d13f-492f-964b-72d34ec1936d.md"
$configName = ".localization-config"
$zhXlfName = "28a0734b-b7aa-416f-b69c-8b1b21830f38.cb1e54aa0bb61f184d790a69febda753cb09c03f.zh-cn.xlf"
$deXlfName = "28a0734b-b7aa-416f-b69c-8b1b21830f38.cb1e54aa0bb61f184d790a69febda753cb09c03f.de-de.xlf"

# ---------------------------------------------------------------------------
# 2. "zh-cn" sheet: the handback report now fills in the "Latest Target
#    File" (E) / "Latest Handback File" (F) columns for both data rows.
#    The handback datetime (G) stays unreported for zh-cn in this run.
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("B2").Value = $newStatus
$wsZh.Range("B3").Value = $newStatus

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $urlMd, "", "", $mdName)
$wsZh.Hyperlinks.Add($wsZh.Range("C2"), $urlZhXlf, "", "", $zhXlfName)
$wsZh.Hyperlinks.Add($wsZh.Range("E2"), $urlMd, "", "", $mdName)
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), $urlZhXlf, "", "", $zhXlfName)
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $urlFfff, "", "", $ffffName)
$wsZh.Hyperlinks.Add($wsZh.Range("C3"), $urlZhXlf, "", "", $zhXlfName)
$wsZh.Hyperlinks.Add($wsZh.Range("E3"), $urlMd, "", "", $mdName)
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), $urlZhXlf, "", "", $zhXlfName)
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), $urlConfig, "", "", $configName)

# ---------------------------------------------------------------------------
# 3. "de-de" sheet: same E/F backfill, plus the handback actually completed
#    so G2/G3 get the real handback timestamp.
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("B2").Value = $newStatus
$wsDe.Range("B3").Value = $newStatus
$wsDe.Range("G2").Value = "2016-03-10 13:04:57"
$wsDe.Range("G3").Value = "2016-03-10 13:04:57"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $urlMd, "", "", $mdName)
$wsDe.Hyperlinks.Add($wsDe.Range("C2"), $urlDeXlf, "", "", $deXlfName)
$wsDe.Hyperlinks.Add($wsDe.Range("E2"), $urlMd, "", "", $mdName)
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), $urlDeXlf, "", "", $deXlfName)
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $urlFfff, "", "", $ffffName)
$wsDe.Hyperlinks.Add($wsDe.Range("C3"), $urlDeXlf, "", "", $deXlfName)
$wsDe.Hyperlinks.Add($wsDe.Range("E3"), $urlMd, "", "", $mdName)
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), $urlDeXlf, "", "", $deXlfName)
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), $urlConfig, "", "", $configName)
